$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update utterance text in column B (rows 1-13): placeholders extracted
# to {PromoName}/{ShowName}/{StartDate}/{EndDate}/{Length}/{AirDate}. ---
$ws.Range("B1").Value  = "When did {PromoName} run?"
$ws.Range("B2").Value  = "Where did {PromoName} run?"
$ws.Range("B3").Value  = "When and where did {PromoName} run?"
$ws.Range("B4").Value  = "What {ShowName} promos are available to run right now?"
$ws.Range("B5").Value  = "Is the {PromoName} for {ShowName} available to run now?"
$ws.Range("B6").Value  = "What promos aired for {ShowName} from {StartDate} to {EndDate}?"
$ws.Range("B7").Value  = "How much have we spent on music for {ShowName}?"
$ws.Range("B8").Value  = "When was the last time {PromoName} aired?"
$ws.Range("B10").Value = "What ran on air on {AirDate}?"
$ws.Range("B11").Value = "Show me all the {Length} second promos available to run for {ShowName}."
$ws.Range("B12").Value = "Give me all promo airings from last night for {ShowName}."
$ws.Range("B13").Value = "Give me all the promos that aired during {ShowName}."

# Rows 1-3 lost their "Notes" column C text (the shared promo-code note was removed).
$ws.Range("C1").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()

# --- Row 13 gets a bottom border under it (closing off the DigitalPromoIntent block). ---
$hdr = $ws.Range("A13:E13")
$hdr.Borders.Item(9).LineStyle = 1
$hdr.Borders.Item(9).Weight = 2

# --- Sheet view: active selection moved. ---
[void]$ws.Range("B14").Select()

# --- Column B widened to match column E. ---
$ws.Columns.Item(2).ColumnWidth = 87.33
